# Adds two new client records to the "Base Clientes carga manual" sheet,
# growing the "Tabla1" Excel table from A1:B19 to A1:B21.
#
#   Row 20 -> Num_Distri 500258    / Name_Distri "AGROINSUMOS DEL CENTRO"
#   Row 21 -> Num_Distri 10268402  / Name_Distri "MAS AGROQUIMICOS Y SEMILLAS"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Tabla1")

# --- new row 20 ---------------------------------------------------------
$lo.ListRows.Add() | Out-Null
$ws.Range("A20").Value = 500258
$ws.Range("B20").Value = "AGROINSUMOS DEL CENTRO"

# --- new row 21 ---------------------------------------------------------
$lo.ListRows.Add() | Out-Null
$ws.Range("A21").Value = 10268402
$ws.Range("B21").Value = "MAS AGROQUIMICOS Y SEMILLAS"

# Leave the newly entered row selected, matching the editor's end state.
$ws.Range("A21:B21").Select()
